# Update the "Run 0" (column B) predicted values on the "Station 2 best results" sheet
# to match the re-run lasso (beta = 1.0) model results for SPI-6 at Station 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new "Run 0" value (column B), rows 2..138
$run0Updates = @{
    2 = 1.549176759265603
    3 = 1.207701839945585
    4 = 0.4059069232809962
    5 = 0.6197262507007525
    6 = -0.5813124571115338
    7 = -0.5947090340710612
    8 = 0.4071683009710324
    9 = 0.2497451395804539
    10 = 0.1071697711574002
    11 = -0.2603564670779595
    12 = 0.2689724442646475
    13 = -0.01062751331263967
    14 = -1.103099120471017
    15 = 0.08401561451310291
    16 = 0.5363350708952966
    17 = 0.476917956832671
    18 = -0.4287162725379494
    19 = 0.001833135909045547
    20 = -0.1398423688360324
    21 = -1.340808767415582
    22 = -1.705854326059518
    23 = -2.039448963493358
    24 = -1.21213279867102
    25 = -1.767482356693216
    26 = -1.342228527536268
    27 = -1.553642930927492
    28 = -1.283857473917002
    29 = -0.009226763556054784
    30 = 0.5052053363459291
    31 = 0.209730297568416
    32 = -0.196899503026408
    33 = 0.3080085646509561
    34 = 0.2584003237712873
    35 = -0.5235714759668169
    36 = -1.293065712645356
    37 = -1.132812276383221
    38 = -0.8608922271349226
    39 = -0.4903253230192564
    40 = -0.379026357962915
    41 = -0.08774943179050432
    42 = -0.1421488945607436
    43 = 0.598658900713158
    44 = 1.687052793765762
    45 = 1.517051733533639
    46 = 1.770712758282326
    47 = 1.693682151837131
    48 = 1.65755313489594
    49 = 1.416136097010809
    50 = 0.1308784473244288
    51 = -0.4480277561205335
    52 = -1.213586179045625
    53 = -0.7484826939065479
    54 = 0.4689450194640229
    55 = 0.8830930660717636
    56 = 1.450643640291943
    57 = 1.595517911349871
    58 = 1.508885076445824
    59 = 1.253709079241698
    60 = 1.153927249036513
    61 = 0.9393366283385978
    62 = 0.2280621403889982
    63 = -0.0312975496032899
    64 = 1.608597886430497
    65 = 1.671621622421648
    66 = 1.727304142887148
    67 = 1.617420229603066
    68 = 1.481691139380476
    69 = 1.646217857679755
    70 = 0.2936229183003242
    71 = 0.6920210224040524
    72 = 0.4216827874561695
    73 = 0.4619377542697602
    74 = 0.6649312871298643
    75 = 0.3711518852719478
    76 = 0.644548530877358
    77 = -0.1509016599632591
    78 = -0.2866072447831352
    79 = 0.6852530318361352
    80 = 0.8157807494219389
    81 = 0.9163382400087725
    82 = 0.1155554131639938
    83 = 0.8444865946937131
    84 = 0.3441985328533933
    85 = -0.7186477712891075
    86 = -1.060943960013575
    87 = -1.492114973436288
    88 = -0.8438698256050454
    89 = -1.46133298002621
    90 = -0.2988994650394377
    91 = -0.02601499975193322
    92 = 0.08947352257780494
    93 = 0.4136285079104567
    94 = 0.5074988852624884
    95 = 0.2134823227425671
    96 = -0.5543124604684329
    97 = -0.4176506659591719
    98 = -0.5392685960225498
    99 = -1.076527023850668
    100 = -0.7999342194309114
    101 = -0.5922446148363774
    102 = -1.157779813179636
    103 = -1.774442223448738
    104 = -2.050775339949804
    105 = -1.458410267496494
    106 = -1.582624969959709
    107 = -0.7645501676742894
    108 = 0.5882406996092652
    109 = 0.7714697408668987
    110 = 1.088764011160676
    111 = 1.671401623200452
    112 = 1.879156717969135
    113 = 1.561988397409557
    114 = 1.259359498391955
    115 = 1.444224059953203
    116 = 1.671590149378062
    117 = 1.393669855896711
    118 = 1.019686208665958
    119 = 1.084217277473481
    120 = 1.886115504459864
    121 = 1.739878254928009
    122 = 1.563125637165057
    123 = 1.221440408119267
    124 = 1.376733040570663
    125 = 1.310278859619054
    126 = -0.6969943775808805
    127 = 0.5901801276701386
    128 = 0.3103700477574086
    129 = 0.8819931068323763
    130 = 0.5209240763436908
    131 = 1.042364817973388
    132 = 1.383269251727831
    133 = 0.5782876576698686
    134 = 0.8608069910223293
    135 = 0.6302438694099107
    136 = 0.3676336302610863
    137 = -0.08975336647095378
    138 = 0.2181652479401538
}

foreach ($row in $run0Updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $run0Updates[$row]
}

Write-Host "Updated $($run0Updates.Count) values in column B (Run 0)"
